# The deck originally ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral"      (used by the slide master / deck)
#   ppt/theme/theme2.xml -> "Office Theme"  (used by the notes master)
#
# The authored change swaps the two themes' contents in place, so that
# ppt/theme/theme1.xml becomes the default "Office Theme" colours (the
# theme actually driving the visible slides/master) while the "Integral"
# colours move to ppt/theme/theme2.xml.
#
# The only real content difference between the two theme parts is their
# 12-slot colour scheme (fonts/format scheme are identical), so we apply
# that swap through the presentation's live theme colour scheme, which is
# the color set that actually renders on the slide master / slides.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Index order (MsoThemeColorSchemeIndex): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
